$p = $ppt.ActivePresentation

function Update-DateShape($shape) {
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "25/05/2021") {
            $tr.Text = "26/05/2021"
        }
    }
}

# Update the slide master's date placeholder
$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    Update-DateShape $m.Shapes.Item($i)
}

# Update every slide layout's date placeholder
for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $layout = $m.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape $layout.Shapes.Item($i)
    }
}
